$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 10 de Julio de 2020 a las 00:41"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("B4").Value = 3214159
$ws.Range("C4").Value = 55227
$ws.Range("D4").Value = 1422416
$ws.Range("E4").Value = 1656063
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 818
$ws.Range("H4").Value = 135680

# Row 8: Peru -> Peru
$ws.Range("B8").Value = 316448
$ws.Range("C8").Value = 3537
$ws.Range("D8").Value = 207802
$ws.Range("E8").Value = 97332
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 181
$ws.Range("H8").Value = 11314

# Row 19: Alemania -> Alemania
$ws.Range("B19").Value = 199198
$ws.Range("C19").Value = 433
$ws.Range("D19").Value = 183600
$ws.Range("E19").Value = 6473
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 10
$ws.Range("H19").Value = 9125

# Row 22: Colombia -> Colombia
$ws.Range("B22").Value = 133973
$ws.Range("C22").Value = 5335
$ws.Range("D22").Value = 56272
$ws.Range("E22").Value = 72987
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 187
$ws.Range("H22").Value = 4714

# Row 52: Armenia -> Nigeria
$ws.Range("A52").Value = "Nigeria"
$ws.Range("B52").Value = 30748
$ws.Range("C52").Value = 499
$ws.Range("D52").Value = 12546
$ws.Range("E52").Value = 17513
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 5
$ws.Range("H52").Value = 689

# Row 53: Nigeria -> Armenia
$ws.Range("A53").Value = "Armenia"
$ws.Range("B53").Value = 30346
$ws.Range("C53").Value = 526
$ws.Range("D53").Value = 18000
$ws.Range("E53").Value = 11811
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 14
$ws.Range("H53").Value = 535

# Row 54: Honduras -> Guatemala
$ws.Range("A54").Value = "Guatemala"
$ws.Range("B54").Value = 26658
$ws.Range("C54").Value = 1247
$ws.Range("D54").Value = 3797
$ws.Range("E54").Value = 21769
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 39
$ws.Range("H54").Value = 1092

# Row 55: Irlanda -> Honduras
$ws.Range("A55").Value = "Honduras"
$ws.Range("B55").Value = 25978
$ws.Range("C55").Value = 550
$ws.Range("D55").Value = 2721
$ws.Range("E55").Value = 22563
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 17
$ws.Range("H55").Value = 694

# Row 56: Guatemala -> Irlanda
$ws.Range("A56").Value = "Irlanda"
$ws.Range("B56").Value = 25565
$ws.Range("C56").Value = 23
$ws.Range("D56").Value = 23364
$ws.Range("E56").Value = 458
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 5
$ws.Range("H56").Value = 1743

# Row 57: Ghana -> Ghana
$ws.Range("B57").Value = 23463
$ws.Range("C57").Value = 641
$ws.Range("D57").Value = 18622
$ws.Range("E57").Value = 4712
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 129

# Row 59: Japon -> Japon
$ws.Range("B59").Value = 20371
$ws.Range("C59").Value = 197
$ws.Range("D59").Value = 17466
$ws.Range("E59").Value = 1924
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 981

# Row 68: Dinamarca -> Chequia
$ws.Range("A68").Value = "Chequia"
$ws.Range("B68").Value = 12919
$ws.Range("C68").Value = 105
$ws.Range("D68").Value = 8128
$ws.Range("E68").Value = 4439
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 352

# Row 69: Chequia -> Dinamarca
$ws.Range("A69").Value = "Dinamarca"
$ws.Range("B69").Value = 12916
$ws.Range("C69").Value = 16
$ws.Range("D69").Value = 12045
$ws.Range("E69").Value = 262
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 609

# Row 71: Uzbekistan -> Uzbekistan
$ws.Range("B71").Value = 11564
$ws.Range("C71").Value = 472
$ws.Range("D71").Value = 7287
$ws.Range("E71").Value = 4226
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 6
$ws.Range("H71").Value = 51

# Row 85: Haiti -> Bulgaria
$ws.Range("A85").Value = "Bulgaria"
$ws.Range("B85").Value = 6672
$ws.Range("C85").Value = 330
$ws.Range("D85").Value = 3229
$ws.Range("E85").Value = 3181
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 3
$ws.Range("H85").Value = 262

# Row 86: Costa Rica -> Haiti
$ws.Range("A86").Value = "Haiti"
$ws.Range("B86").Value = 6486
$ws.Range("C86").Value = 54
$ws.Range("D86").Value = 2181
$ws.Range("E86").Value = 4182
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 6
$ws.Range("H86").Value = 123

# Row 87: Tayikistan -> Costa Rica
$ws.Range("A87").Value = "Costa Rica"
$ws.Range("B87").Value = 6485
$ws.Range("C87").Value = 649
$ws.Range("D87").Value = 2023
$ws.Range("E87").Value = 4437
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 25

# Row 88: Bulgaria -> Tayikistan
$ws.Range("A88").Value = "Tayikistan"
$ws.Range("B88").Value = 6410
$ws.Range("C88").Value = 46
$ws.Range("D88").Value = 5067
$ws.Range("E88").Value = 1289
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 54

# Row 90: Gabon -> Guinea
$ws.Range("A90").Value = "Guinea"
$ws.Range("B90").Value = 5881
$ws.Range("C90").Value = 184
$ws.Range("D90").Value = 4672
$ws.Range("E90").Value = 1173
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 2
$ws.Range("H90").Value = 36

# Row 91: Guinea -> Gabon
$ws.Range("A91").Value = "Gabon"
$ws.Range("B91").Value = 5871
$ws.Range("C91").Value = 0
$ws.Range("D91").Value = 2682
$ws.Range("E91").Value = 3143
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 46

# Row 94: Mauritania -> Mauritania
$ws.Range("B94").Value = 5126
$ws.Range("C94").Value = 39
$ws.Range("D94").Value = 2026
$ws.Range("E94").Value = 2956
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 5
$ws.Range("H94").Value = 144

# Row 98: Republica de Africa Central -> Republica de Africa Central
$ws.Range("B98").Value = 4200
$ws.Range("C98").Value = 91
$ws.Range("D98").Value = 1142
$ws.Range("E98").Value = 3006
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 52

# Row 114: Libano -> Libano
$ws.Range("B114").Value = 2011
$ws.Range("C114").Value = 65
$ws.Range("D114").Value = 1368
$ws.Range("E114").Value = 607
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 36

# Row 184: Seychelles -> Lesoto
$ws.Range("A184").Value = "Lesoto"
$ws.Range("B184").Value = 91
$ws.Range("C184").Value = 0
$ws.Range("D184").Value = 11
$ws.Range("E184").Value = 80
$ws.Range("F184").Value = 0
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 0

# Row 185: Lesoto -> Seychelles
$ws.Range("A185").Value = "Seychelles"
$ws.Range("B185").Value = 91
$ws.Range("C185").Value = 0
$ws.Range("D185").Value = 11
$ws.Range("E185").Value = 80
$ws.Range("F185").Value = 0
$ws.Range("G185").Value = 0
$ws.Range("H185").Value = 0
